# Lesson 14 wordlist update: restore the correct English/Japanese
# pairings for rows that had drifted out of sync with the lesson
# sections (also fixes the Lesson 13 section placeholder tail).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @(2, '(my) older brother', '兄|あに'),
    @(3, 'landlord; landlady', '大家さん|おおやさん'),
    @(4, 'return (as a token of gratitude)', 'お返し|おかえし'),
    @(5, '(your/his) wife', '奥さん|おくさん'),
    @(6, 'uncle; middle-aged man', 'おじさん'),
    @(7, 'aunt; middle-aged woman', 'おばさん'),
    @(8, 'tumbler; glass', 'グラス'),
    @(9, 'Christmas', 'クリスマス'),
    @(10, '(your/his) husband', 'ご主人|ごしゅじん'),
    @(11, 'plate; dish', '皿|さら'),
    @(12, 'time', '時間|じかん'),
    @(13, 'ticket', 'チケット'),
    @(14, 'chocolate', 'チョコレート'),
    @(15, 'sweat shirt', 'トレーナー'),
    @(16, 'stuffed animal', 'ぬいぐるみ'),
    @(17, 'necktie', 'ネクタイ'),
    @(18, 'St. Valentine''s Day', 'バレンタインデー'),
    @(19, 'camcorder', 'ビデオカメラ'),
    @(20, 'married couple; husband and wife', '夫婦|ふうふ'),
    @(21, 'White Day', 'ホワイトデー'),
    @(22, 'winter scarf', 'マフラー'),
    @(23, 'comic book', '漫画|まんが'),
    @(24, 'multistory apartment building; condo', 'マンション'),
    @(25, 'mandarin orange', 'みかん'),
    @(26, 'everyone; all of you', '皆さん|みなさん'),
    @(27, 'ring', '指輪|ゆびわ'),
    @(28, 'radio', 'ラジオ'),
    @(29, 'parents', '両親|りょうしん'),
    @(30, 'résumé', '履歴書|りれきしょ'),
    @(31, 'to want', '欲しい|ほしい'),
    @(32, 'stingy; cheap', 'けち（な）'),
    @(33, 'to send', '送る|おくる'),
    @(34, 'to look good (on somebody)', '似合う|にあう'),
    @(35, 'to give up', 'あきらめる'),
    @(36, 'to give (to others)', 'あげる'),
    @(37, 'to give (me)', 'くれる'),
    @(38, 'to come into existence; to be made', 'できる'),
    @(39, 'to consult', '相談する|そうだんする'),
    @(40, 'to propose marriage', 'プロポーズする'),
    @(41, 'same', '同じ|おなじ'),
    @(42, 'Mr./Ms....(casual)', '～君|～くん'),
    @(43, '...like this; this kind of...', 'こんな～'),
    @(44, '[makes a noun plural]', '～たち'),
    @(45, 'we', '私たち|わたしたち'),
    @(46, 'exactly', 'ちょうど'),
    @(47, 'what should one do', 'どうしたらいい'),
    @(48, 'well', 'よく'),
    @(55, 'he; boyfriend', '彼|かれ'),
    @(56, 'she; girlfriend', '彼女|かのじょ'),
    @(57, 'they', '彼ら|かれら'),
    @(58, 'age; era', '時代|じだい'),
    @(59, 'electricity fee', '電気代|でんきだい'),
    @(60, '90''s', '九十年代|きゅうじゅうねんだい'),
    @(61, 'in one''s teens', '十代|じゅうだい'),
    @(62, 'instead', '代わりに|かわりに'),
    @(63, 'foreign students', '留学生|りゅうがくせい'),
    @(64, 'to study abroad', '留学する|りゅうがくする'),
    @(65, 'absence from home', '留守|るす'),
    @(66, 'family', '家族|かぞく'),
    @(67, 'race', '民族|みんぞく'),
    @(68, 'aquarium', '水族館|すいぞくかん'),
    @(69, 'member of royalty', '王族|おうぞく'),
    @(70, 'father', '父親|ちちおや'),
    @(71, 'kind', '親切な|しんせつな'),
    @(72, 'best friend', '親友|しんゆう'),
    @(73, 'parents', '両親|りょうしん'),
    @(74, 'intimate', '親しい|したしい'),
    @(75, 'mother', '母親|ははおや'),
    @(76, 'to cut', '切る|きる'),
    @(77, 'ticket', '切符|きっぷ'),
    @(78, 'postage stamp', '切手|きって'),
    @(79, 'precious', '大切な|たいせつな'),
    @(80, 'English', '英語|えいご'),
    @(81, 'England', '英国|えいこく'),
    @(82, 'English conversation', '英会話|えいかいわ'),
    @(83, 'hero', '英雄|えいゆう'),
    @(84, 'shop', '店|みせ'),
    @(85, 'store attendant', '店員|てんいん'),
    @(86, 'stall', '売店|ばいてん'),
    @(87, 'book store', '書店|しょてん'),
    @(88, 'last year', '去年|きょねん'),
    @(89, 'the past', '過去|かこ'),
    @(90, 'to leave', '去る|さる'),
    @(91, 'to erase', '消去する|しょうきょする'),
    @(104, 'music', '音楽|おんがく'),
    @(105, 'pronunciation', '発音|はつおん'),
    @(106, 'sound', '音|おと'),
    @(107, 'real intention', '本音|ほんね'),
    @(108, 'fun', '楽しい|たのしい'),
    @(109, 'musical instrument', '楽器|がっき'),
    @(110, 'easy; comfortable', '楽な|らくな'),
    @(111, 'doctor', '医者|いしゃ'),
    @(112, 'dentist', '歯医者|はいしゃ'),
    @(113, 'medical science', '医学|いがく'),
    @(114, 'clinic', '医院|いいん'),
    @(115, 'scholar', '学者|がくしゃ'),
    @(116, 'reader', '読者|どくしゃ'),
    @(117, 'young people', '若者|わかもの')
)

foreach ($row in $rows) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
}
